# Weekly price-update commit: a new observation row is inserted at row 169
# (pushing the existing rows 169-253 down to 170-254) with fresh data for
# "Provincia del Elquí".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 169, shifting rows 169:253 down to 170:254.
$ws.Rows("169:169").Insert()

# Populate the newly inserted row 169 with the new record.
$ws.Range("A169").Value = 9
$ws.Range("B169").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C169").Value = "Metropolitana"
$ws.Range("D169").Value = 44784
$ws.Range("E169").Value = 13
$ws.Range("F169").Value = 100112026
$ws.Range("G169").Value = "Haba"
$ws.Range("H169").Value = "Sin especificar"
$ws.Range("I169").Value = "Primera"
$ws.Range("J169").Value = 62
$ws.Range("K169").Value = 13000
$ws.Range("L169").Value = 14000
$ws.Range("M169").Value = 13645
$ws.Range("N169").Value = "`$/saco 25 kilos"
$ws.Range("O169").Value = "Provincia del Elquí"
$ws.Range("P169").Value = 546
$ws.Range("Q169").Value = 25
$ws.Range("R169").Value = "Hortaliza"
